$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column-by-column to reproduce the original authoring order
# (this determines the order new entries are appended to the shared
# strings table).

# Column A - Request ID
$ws.Range("A3").Value = "IR003"
$ws.Range("A4").Value = "IR003"

# Column B - Applicant Name
$ws.Range("B4").Value = "Sujata"
$ws.Range("B5").Value = "Sujata"

# Column H - Country
$ws.Range("H3").Value = "IND"

# Column E - Fathers Name
$ws.Range("E3").Value = "A"
$ws.Range("E4").Value = "B"
$ws.Range("E5").Value = "C"

# Column D - Address
$ws.Range("D3").Value = "PUNE"
$ws.Range("D4").Value = "MUMBAI"
$ws.Range("D5").Value = "DELHI"

# Match final selection shown in the diff
$ws.Range("D5").Select()
